$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in the missing "Taille" (height) values for a few players
$ws.Range("E19").Value = "1m84"
$ws.Range("E24").Value = "1m89"
$ws.Range("E26").Value = "1m74"

# Move the active selection to D31 (matches the saved cursor position)
$ws.Range("D31").Select()
